# DeleteEmployee-BackEnd-Admin.xlsx update
# Adds two new columns of data to row 5 ("an exception was thrown" / "fail")
# and moves the sheet's scroll/selection so F5 is the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cells added to row 5.
$ws.Range("E5").Value = "an exception was thrown"
$ws.Range("F5").Value = "fail"

# Match the wrap-text / top-vertical-aligned formatting used by the rest
# of the table (same look as C5/D5, style index 1 in the original file).
$ws.Range("E5:F5").WrapText = $true
$ws.Range("E5:F5").VerticalAlignment = -4160

# Update the view: scroll so column E is the left-most visible column
# (topLeftCell goes from A3 to E3) and make F5 the active/selected cell
# (was D7).
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 3
$ws.Range("F5").Select()
